# Generate Report for Handback
#
# The localization-status report is regenerated after a successful
# handback: the target language is back in sync with en-US, so the
# "Ready for handoff" status becomes "Handed back: in sync with en-US",
# the "Latest Handback DateTime" timestamps are refreshed, and the
# stale "version mismatch" Error Detail message is cleared since the
# handback file is now current.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns -------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# --- zh-cn detail sheet ----------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-27 12:46:48"
$zhcn.Range("P2").Value = ""

# --- de-de detail sheet -----------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-27 12:46:54"
$dede.Range("P2").Value = ""

# --- Column widths: the Status column is now wider (longer text) and
#     the now-empty Error Detail column is narrower. The host rounds
#     ColumnWidth to the nearest 1/6 character, so these inputs are the
#     closest achievable approximations of the recorded widths.
$statusWidth = 29.166666666666668
$errorWidth = 12.833333333333334

$overview.Range("E1").EntireColumn.ColumnWidth = $statusWidth
$overview.Range("F1").EntireColumn.ColumnWidth = $statusWidth

$zhcn.Range("C1").EntireColumn.ColumnWidth = $statusWidth
$zhcn.Range("P1").EntireColumn.ColumnWidth = $errorWidth

$dede.Range("C1").EntireColumn.ColumnWidth = $statusWidth
$dede.Range("P1").EntireColumn.ColumnWidth = $errorWidth
